# Fixed battery and signal DSL issues
#
# 1. Row 2 (G2): StartPage/WebServer config step had FullScreen value="0" -
#    flip it to value="1".
# 2. Rows 16-18 (VT298_0003 / VT298_0006 / VT298_0012): the H-column
#    validate4 block used a Screenshot validation; swap it for an
#    isIconDisplayed validation against the signal view icon.
# 3. Row 18 (VT298_0012), G column: drop the stray TakeScreenshot(VT298_0012)
#    step that isn't needed any more now that validate4 no longer screenshots.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# --- 1. G2: FullScreen 0 -> 1 --------------------------------------------
$g2 = $ws.Range("G2").Text
$g2 = $g2.Replace('FullScreen value="0"', 'FullScreen value="1"')
$ws.Range("G2").Value = $g2

# --- 2. G18: remove the extraneous TakeScreenshot(VT298_0012); line -------
$g18 = $ws.Range("G18").Text
$g18 = $g18.Replace("TakeScreenshot(VT298_0012);`n", "")
$ws.Range("G18").Value = $g18

# --- 3. H16 / H17 / H18: Screenshot validation -> isIconDisplayed ---------
$h16 = $ws.Range("H16").Text
$h16 = $h16.Replace("validate_Screenshot=VT298_0003", "validate_isIconDisplayed=signalview_xpath,true")
$ws.Range("H16").Value = $h16

$h17 = $ws.Range("H17").Text
$h17 = $h17.Replace("validate_Screenshot=VT298_0006", "validate_isIconDisplayed=signalview_xpath,true")
$ws.Range("H17").Value = $h17

$h18 = $ws.Range("H18").Text
$h18 = $h18.Replace("validate_Screenshot=VT298_0012", "validate_isIconDisplayed=signalview_xpath,true")
$ws.Range("H18").Value = $h18

# --- 4. Leave the cursor parked on D1, matching the saved selection -------
$ws.Range("D1").Select()
